# Apply updated cryptocurrency price/volume data per Feb 18 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.653.78"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "1.690.65"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  +0.53%  "

$ws.Range("D5").Value = "'315.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("E6").Value = "  +0.54%  "

$ws.Range("D7").Value = "'0.3940"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.4062"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").Value = "'1.001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("D11").Value = "'52.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.18%  "

$ws.Range("D12").Value = "'0.08811"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").Value = "'7.226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "'23.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("D15").Value = "'8.059"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.16%  "

$ws.Range("D17").Value = "1.697.87"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "'100.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").Value = "'0.07011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'19.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").Value = "'7.013"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.77%  "

$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").Value = "'14.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").Value = "24.660.56"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").Value = "'3.261"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.02%  "

$ws.Range("D26").Value = "'2.364"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.67%  "

$ws.Range("D27").Value = "'22.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.37%  "

$ws.Range("E28").Value = "  +2.42%  "

$ws.Range("D29").Value = "'135.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.00%  "

$ws.Range("D30").Value = "'5.194"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").Value = "'7.593"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.73%  "

$ws.Range("D32").Value = "1.884.91"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("E33").Value = "  -1.70%  "

$ws.Range("D34").Value = "'1.054"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.69%  "

$ws.Range("D35").Value = "'7.174"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("D36").Value = "'11.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").Value = "'1.906"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("E39").Value = "  -2.66%  "

$ws.Range("D40").Value = "'0.09189"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.48%  "

$ws.Range("D41").Value = "'0.02719"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.91%  "

$ws.Range("D42").Value = "'1.458"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").Value = "'0.7610"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "'16.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.70%  "

$ws.Range("D45").Value = "'0.7135"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("D46").Value = "'2.570"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.01%  "

$ws.Range("D47").Value = "'4.209"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").Value = "'1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("D49").Value = "'1.317"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.47%  "

$ws.Range("D50").Value = "'139.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("E51").Value = "  +0.01%  "
